# "Working extra client side validation" -
# add a new "Add"/"Hinzufügen" resource row and a new "Nationality"/
# "Nationalität" resource row to both the "en" and "de" sheets.

$wb = $excel.ActiveWorkbook
$wsEn = $wb.Worksheets.Item("en")
$wsDe = $wb.Worksheets.Item("de")

# New shared strings must be created in this order: Add, Hinzufügen,
# Nationality, Nationalität - so row 58 (both sheets) is written fully
# before row 59 is touched.

# --- row 58: carries the same wrap/vertical-center format as the rest of
# the table, so clone it from an existing formatted row instead of setting
# alignment properties piecemeal (which would mint extra, unused styles).
$wsEn.Range("A2:B2").Copy()
$wsEn.Range("A58:B58").PasteSpecial(-4122)
$wsDe.Range("A2:B2").Copy()
$wsDe.Range("A58:B58").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsEn.Cells.Item(58, 1).Value = "Add"
$wsEn.Cells.Item(58, 2).Value = "Add"

$wsDe.Cells.Item(58, 1).Value = "Add"
$wsDe.Cells.Item(58, 2).Value = "Hinzufügen"

# --- row 59: plain/default formatting (no style), as in the source sheet ---
$wsEn.Cells.Item(59, 1).Value = "Nationality"
$wsEn.Cells.Item(59, 2).Value = "Nationality"

$wsDe.Cells.Item(59, 1).Value = "Nationality"
$wsDe.Cells.Item(59, 2).Value = "Nationalität"

# B59 on the "de" sheet keeps the formatted style, same as column B above it.
$wsDe.Range("B2").Copy()
$wsDe.Range("B59").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsDe.Cells.Item(59, 2).Value = "Nationalität"

# --- refresh the active-cell selection on each sheet to B59 ---
$wsEn.Activate()
$wsEn.Range("B59").Select()

$wsDe.Activate()
$wsDe.Range("B59").Select()

# "en" was the originally-selected tab - leave it active.
$wsEn.Activate()
